# Commit: "Fruta / hortaliza, semanal"
# Insert two new weekly price-report rows right after the existing row 141
# (pushing the old rows 142-191 down to 144-193) and populate the two new
# rows with the new "Camote" (Zapallo) observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 142:143 - everything from the old row 142 onward
# shifts down by two rows (old 142 -> new 144, ..., old 191 -> new 193).
$ws.Rows("142:143").Insert()

# --- New row 142 ---------------------------------------------------------
$ws.Cells.Item(142, 1).Value  = 7
$ws.Cells.Item(142, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(142, 3).Value  = "Ñuble"
$ws.Cells.Item(142, 4).Value  = 44845
$ws.Cells.Item(142, 5).Value  = 16
$ws.Cells.Item(142, 6).Value  = 100112045
$ws.Cells.Item(142, 7).Value  = "Zapallo"
$ws.Cells.Item(142, 8).Value  = "Camote"
$ws.Cells.Item(142, 9).Value  = "1a (guarda)"
$ws.Cells.Item(142, 10).Value = 300
$ws.Cells.Item(142, 11).Value = 900
$ws.Cells.Item(142, 12).Value = 1000
$ws.Cells.Item(142, 13).Value = 950
$ws.Cells.Item(142, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(142, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(142, 16).Value = 950
$ws.Cells.Item(142, 17).Value = 1
$ws.Cells.Item(142, 18).Value = "Hortaliza"

# --- New row 143 ---------------------------------------------------------
$ws.Cells.Item(143, 1).Value  = 7
$ws.Cells.Item(143, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(143, 3).Value  = "Ñuble"
$ws.Cells.Item(143, 4).Value  = 44845
$ws.Cells.Item(143, 5).Value  = 16
$ws.Cells.Item(143, 6).Value  = 100112045
$ws.Cells.Item(143, 7).Value  = "Zapallo"
$ws.Cells.Item(143, 8).Value  = "Camote"
$ws.Cells.Item(143, 9).Value  = "2a (guarda)"
$ws.Cells.Item(143, 10).Value = 200
$ws.Cells.Item(143, 11).Value = 800
$ws.Cells.Item(143, 12).Value = 800
$ws.Cells.Item(143, 13).Value = 800
$ws.Cells.Item(143, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(143, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(143, 16).Value = 800
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = "Hortaliza"
